$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.85
$ws.Range("H2").Value = 4.3
$ws.Range("N2").Value = 4.5
$ws.Range("P2").Value = 2.24
$ws.Range("R2").Value = 1.49
$ws.Range("S2").Value = 2.86
$ws.Range("T2").Value = 1.7
$ws.Range("U2").Value = 2.28
$ws.Range("AA2").Value = 110
$ws.Range("AC2").Value = 9.6
$ws.Range("AH2").Value = 980
$ws.Range("AM2").Value = 100

$ws.Range("G4").Value = 2.06
$ws.Range("H4").Value = 3.45
$ws.Range("I4").Value = 5.6
$ws.Range("J4").Value = 3.65
$ws.Range("Q4").Value = 1.78

$ws.Range("G5").Value = 4.3
$ws.Range("H5").Value = 1.85
$ws.Range("J5").Value = 4.1
$ws.Range("K5").Value = 4.7
$ws.Range("P5").Value = 2.52
$ws.Range("Q5").Value = 1.52

$ws.Range("N6").Value = 3.7
$ws.Range("P6").Value = 1.92
$ws.Range("R6").Value = 1.35
$ws.Range("S6").Value = 3.6
$ws.Range("T6").Value = 1.78
$ws.Range("U6").Value = 2.18

$ws.Range("Y8").Value = 9.4
$ws.Range("Z8").Value = 32
